$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# sheet1
$ws1.Range("F4").Value = 13
$ws1.Range("F6").Value = 20
$ws1.Range("F7").Value = 854
$ws1.Range("F8").Value = 37
$ws1.Range("F9").Value = 6791
$ws1.Range("F12").Value = 138
$ws1.Range("F13").Value = 6425
$ws1.Range("F14").Value = 125
$ws1.Range("F16").Value = 4342
$ws1.Range("F20").Value = 4306
$ws1.Range("F21").Value = 222
$ws1.Range("F22").Value = 229
$ws1.Range("F23").Value = 313
$ws1.Range("F24").Value = 272
$ws1.Range("F25").Value = 226
$ws1.Range("F26").Value = 125
$ws1.Range("F28").Value = 40
$ws1.Range("F31").Value = 68
$ws1.Range("F32").Value = 7853
$ws1.Range("F33").Value = 50
$ws1.Range("F34").Value = 1330
$ws1.Range("F35").Value = 650
$ws1.Range("F38").Value = 984
$ws1.Range("F39").Value = 1564
$ws1.Range("F41").Value = 903
$ws1.Range("F43").Value = 3919
$ws1.Range("F46").Value = 106
$ws1.Range("F47").Value = 31
$ws1.Range("F49").Value = 1080

# sheet3
$ws3.Range("F2").Value = 232

# sheet4
$ws4.Range("F2").Value = 232
$ws4.Range("F7").Value = 13
$ws4.Range("F9").Value = 20
$ws4.Range("F10").Value = 854
$ws4.Range("F11").Value = 37
$ws4.Range("F12").Value = 6791
$ws4.Range("F15").Value = 138
$ws4.Range("F16").Value = 6425
$ws4.Range("F17").Value = 125
$ws4.Range("F19").Value = 4342
$ws4.Range("F22").Value = 4306
$ws4.Range("F23").Value = 222
$ws4.Range("F24").Value = 229
$ws4.Range("F25").Value = 313
$ws4.Range("F26").Value = 272
$ws4.Range("F27").Value = 226
$ws4.Range("F28").Value = 125
$ws4.Range("F31").Value = 68
$ws4.Range("F33").Value = 7853
$ws4.Range("F34").Value = 50
$ws4.Range("F35").Value = 1330
$ws4.Range("F36").Value = 650
$ws4.Range("F38").Value = 984
$ws4.Range("F39").Value = 1564
$ws4.Range("F41").Value = 903
$ws4.Range("F43").Value = 3919
$ws4.Range("F46").Value = 106
$ws4.Range("F48").Value = 1080
